$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full target data (rows 2-11, columns A-T) reflecting updated TPM values
# and a new "MuSCs" sending-cluster block (rows 7-11).
$data = @(
  @("FAPs","Adam2","Itga9","ECs",1,0.3333333333333333,0.1251886666666667,0.375566,0.6155719715657366,0.7060466830097307,3,1,2.621797333333333,7.865392,0.07867217155043885,0.07906089226781998,0.3282193124302222,2.953973811872,0.04842838374866149,0.05582068074148397),
  @("FAPs","Adam2","Itga9","FAPs",1,0.3333333333333333,0.1251886666666667,0.375566,0.6155719715657366,0.7060466830097307,3,1,4.942469,14.827407,0.148308476825081,0.1490412718702539,0.6187411041513333,5.568669937362,0.09129454147912648,0.1052300956355443),
  @("FAPs","Adam2","Itga9","Inflammatory-Mac",1,0.3333333333333333,0.1251886666666667,0.375566,0.6155719715657366,0.7060466830097307,3,1,11.09754033333333,33.292621,0.3330034651388949,0.3346488416844782,1.389286277609556,12.503576498486,0.2049875995737716,0.2362777046443743),
  @("FAPs","Adam2","Itga9","MuSCs",1,0.3333333333333333,0.1251886666666667,0.375566,0.6155719715657366,0.7060466830097307,2,1,0.4915585,0.983117,0.01475017696730553,0.009882038584175128,0.06153755320366667,0.369225319222,0.009079795516707782,0.006977180563731025),
  @("FAPs","Adam2","Itga9","Resolving-Mac",1,0.3333333333333333,0.1251886666666667,0.375566,0.6155719715657366,0.7060466830097307,3,1,14.172235,42.516705,0.4252657095182797,0.4273669555932728,1.774203203336667,15.96782883003,0.2617816512474693,0.3017410214245971),
  @("MuSCs","Adam2","Itga9","ECs",1,0.5,0.078181,0.156362,0.3844280284342634,0.2939533169902694,3,1,2.621797333333333,7.865392,0.07867217155043885,0.07906089226781998,0.2049747373173333,1.229848423904,0.03024378780177735,0.02324021152633602),
  @("MuSCs","Adam2","Itga9","FAPs",1,0.5,0.078181,0.156362,0.3844280284342634,0.2939533169902694,3,1,4.942469,14.827407,0.148308476825081,0.1490412718702539,0.386407168889,2.318443013334,0.05701393534595454,0.04381117623470966),
  @("MuSCs","Adam2","Itga9","Inflammatory-Mac",1,0.5,0.078181,0.156362,0.3844280284342634,0.2939533169902694,3,1,11.09754033333333,33.292621,0.3330034651388949,0.3346488416844782,0.8676168008003333,5.205700804801999,0.1280158655651233,0.09837113704010389),
  @("MuSCs","Adam2","Itga9","MuSCs",1,0.5,0.078181,0.156362,0.3844280284342634,0.2939533169902694,2,1,0.4915585,0.983117,0.01475017696730553,0.009882038584175128,0.0384305350885,0.153722140354,0.005670381450597747,0.002904858020444104),
  @("MuSCs","Adam2","Itga9","Resolving-Mac",1,0.5,0.078181,0.156362,0.3844280284342634,0.2939533169902694,3,1,14.172235,42.516705,0.4252657095182797,0.4273669555932728,1.107999504535,6.647997027210001,0.1634840582708104,0.1256259341686757),
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rowValues = $data[$i]
    for ($c = 0; $c -lt $rowValues.Count; $c++) {
        $ws.Cells.Item($row, $c + 1).Value = $rowValues[$c]
    }
}
